$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column C to fit the new, longer "Observaciones" text
$ws.Columns.Item(3).ColumnWidth = 49.5

# Add the new tracking row (Nro Exposicion 5)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = [datetime]"2012-05-22"
$ws.Range("B6").NumberFormat = "m/d/yy"
$ws.Range("C6").Value = "Exposicion oral /entrega final De Exposicion"
$ws.Range("D6").Value = "Pendiente"

# Move the active selection to C6, matching where the new comment was entered
[void]$ws.Range("C6").Select()
